$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nodes")

$ws.Range("D3").Value = 0.65
$ws.Range("D4").Value = 0.34
$ws.Range("D5").Value = 0.65

$ws.Activate()
$ws.Range("H5").Select()
